# Insert a new weekly record at row 211, shifting the existing rows
# 211-248 down to 212-249 (matching the "semanal" update described in
# the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 211, pushing the rest
# of the table (including row 248) down by one row.
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A211").Value = 4
$ws.Range("B211").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C211").Value = "Los Lagos"
$ws.Range("D211").Value = 44522
$ws.Range("E211").Value = 10
$ws.Range("F211").Value = 100114013
$ws.Range("G211").Value = "Zanahoria"
$ws.Range("H211").Value = "Sin especificar"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 250
$ws.Range("K211").Value = 12000
$ws.Range("L211").Value = 12000
$ws.Range("M211").Value = 12000
$ws.Range("N211").Value = "`$/saco 20 kilos"
$ws.Range("O211").Value = "Región de Ñuble"
$ws.Range("P211").Value = 600
$ws.Range("Q211").Value = 20
$ws.Range("R211").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of the table.
$ws.Range("D211").NumberFormat = $ws.Range("D212").NumberFormat
